$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Renumber Pull # column (A2:A8) 6,7,1,5,2,3,4 -> 1,2,3,4,5,6,7
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

# Add two new rows of data (rows 9 and 10)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "LOCAL"
$ws.Range("C9").Value = "543+00"
$ws.Range("D9").Value = "553+00"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = "2C#4"
$ws.Range("G9").Value = "PK"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "LOCAL"
$ws.Range("C10").Value = "543+00"
$ws.Range("D10").Value = "553+00"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = "2C#4"
$ws.Range("G10").Value = "PK"
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0

$ws.Range("G10").Select()
